$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps its literal text representation
# (values like "0.430", "8.60", "1.10" must not be normalized to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.677.53"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.217.28"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D5").Value = "240.58"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "74.79"
$ws.Range("E7").Value = "  +4.95%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "41.13"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "54.52"
$ws.Range("E12").Value = "  -12.96%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D15").Value = "2.547.80"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "14.69"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").Value = "2.214.78"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "0.801"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "42.506.39"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "70.73"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "5.92"
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D23").Value = "9.81"
$ws.Range("E23").Value = "  -7.87%  "
$ws.Range("D24").Value = "229.44"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  +7.33%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("E28").Value = "  -6.90%  "
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "172.87"
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").Value = "36.49"
$ws.Range("E32").Value = "  +20.97%  "
$ws.Range("D33").Value = "20.27"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "5.28"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("D39").Value = "0.0323"
$ws.Range("E39").Value = "  +8.84%  "
$ws.Range("D40").Value = "12.38"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "60.24"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "8.60"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "0.0987"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "99.17"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.27"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.430"
$ws.Range("E51").Value = "  +18.75%  "
